# Update the "想去人数" (want-to-attend count) column F for specific rows
# on both the "展览" and "全部类型" worksheets, matching the regenerated
# site data output.

$wb = $excel.ActiveWorkbook

# Row number (key) -> new value for column F
$updates = @{
    2  = 26
    8  = 462
    11 = 577
    13 = 302
    15 = 374
    17 = 94
    20 = 50
    22 = 947
    25 = 330
    31 = 222
    34 = 1624
    37 = 160
    40 = 3690
    41 = 429
    43 = 912
    46 = 69
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
